# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows
# on Sheet1, reflecting corrected dialog-act annotations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @(
    @{Row=2;   I="aa"; J="Agree/Accept"},
    @{Row=5;   I="sd"; J="Statement-non-opinion"},
    @{Row=15;  I="sv"; J="Statement-opinion"},
    @{Row=19;  I="aa"; J="Agree/Accept"},
    @{Row=24;  I="b";  J="Acknowledge (Backchannel)"},
    @{Row=31;  I="sd"; J="Statement-non-opinion"},
    @{Row=33;  I="sv"; J="Statement-opinion"},
    @{Row=41;  I="aa"; J="Agree/Accept"},
    @{Row=46;  I="sd"; J="Statement-non-opinion"},
    @{Row=52;  I="sd"; J="Statement-non-opinion"},
    @{Row=73;  I="sd"; J="Statement-non-opinion"},
    @{Row=79;  I="sd"; J="Statement-non-opinion"},
    @{Row=90;  I="aa"; J="Agree/Accept"},
    @{Row=96;  I="sd"; J="Statement-non-opinion"},
    @{Row=97;  I="sd"; J="Statement-non-opinion"},
    @{Row=98;  I="sd"; J="Statement-non-opinion"},
    @{Row=101; I="aa"; J="Agree/Accept"},
    @{Row=102; I="ba"; J="Appreciation"},
    @{Row=103; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=105; I="%";  J="Uninterpretable"},
    @{Row=129; I="ba"; J="Appreciation"},
    @{Row=166; I="ba"; J="Appreciation"},
    @{Row=167; I="sv"; J="Statement-opinion"},
    @{Row=171; I="sv"; J="Statement-opinion"},
    @{Row=174; I="sv"; J="Statement-opinion"},
    @{Row=183; I="sv"; J="Statement-opinion"},
    @{Row=184; I="sd"; J="Statement-non-opinion"},
    @{Row=202; I="%";  J="Uninterpretable"},
    @{Row=209; I="aa"; J="Agree/Accept"},
    @{Row=218; I="sd"; J="Statement-non-opinion"},
    @{Row=220; I="sd"; J="Statement-non-opinion"},
    @{Row=222; I="%";  J="Uninterpretable"},
    @{Row=225; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=234; I="aa"; J="Agree/Accept"},
    @{Row=242; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=248; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=249; I="aa"; J="Agree/Accept"},
    @{Row=252; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=261; I="sd"; J="Statement-non-opinion"},
    @{Row=268; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=270; I="sv"; J="Statement-opinion"},
    @{Row=287; I="aa"; J="Agree/Accept"},
    @{Row=299; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=301; I="b";  J="Acknowledge (Backchannel)"},
    @{Row=309; I="sd"; J="Statement-non-opinion"},
    @{Row=314; I="sd"; J="Statement-non-opinion"},
    @{Row=320; I="sd"; J="Statement-non-opinion"},
    @{Row=324; I="sd"; J="Statement-non-opinion"},
    @{Row=334; I="b";  J="Acknowledge (Backchannel)"}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}
